# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# to reflect the newly generated report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - Latest HO Xliff Generate Date (also shared with de-de Correspond Handoff Datetime)
$wsOverview.Range("G2").Value = "2016-08-19 09:02:03"

# zh-cn sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-19 09:01:55"
$wsZhCn.Range("K2").Value = "2016-08-19 09:02:28"

# de-de sheet - Correspond Handoff Datetime (stays in sync with Overview G2) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-19 09:02:03"
$wsDeDe.Range("K2").Value = "2016-08-19 09:02:34"
